$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "585.50"); force text
# formatting before assignment so Excel does not silently coerce them to
# numbers and drop significant trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.984.02'
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.050.66'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.50'
$ws.Range("E5").Value = '  -0.86%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.55'
$ws.Range("E6").Value = '  -2.10%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.535'
$ws.Range("E8").Value = '  -1.99%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.051.02'
$ws.Range("E9").Value = '  -0.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("E10").Value = '  -2.99%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.85'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("E12").Value = '  -2.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000234'
$ws.Range("E13").Value = '  -3.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.15'
$ws.Range("E14").Value = '  -4.02%  '

$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.556.83'
$ws.Range("E16").Value = '  -0.21%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.11'
$ws.Range("E17").Value = '  -1.47%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.023.38'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.054.47'
$ws.Range("E19").Value = '  -0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.96'
$ws.Range("E20").Value = '  -0.40%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.25'
$ws.Range("E21").Value = '  -3.21%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.704'
$ws.Range("E22").Value = '  -2.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.48'
$ws.Range("E23").Value = '  -1.36%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("E24").Value = '  -0.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.92'
$ws.Range("E25").Value = '  +0.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.62'
$ws.Range("E26").Value = '  -2.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.46'
$ws.Range("E27").Value = '  +4.08%  '

$ws.Range("E28").Value = '  +0.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("E29").Value = '  -0.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.66'
$ws.Range("E30").Value = '  -1.07%  '

$ws.Range("E31").Value = '  +0.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.19'
$ws.Range("E32").Value = '  -0.77%  '

$ws.Range("E33").Value = '  +1.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").Value = '  -3.19%  '

$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0816'
$ws.Range("E36").Value = '  -4.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.23'
$ws.Range("E37").Value = '  -4.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.89'
$ws.Range("E38").Value = '  -3.79%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.20'
$ws.Range("E39").Value = '  -1.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.21'
$ws.Range("E40").Value = '  -1.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.34'
$ws.Range("E41").Value = '  -0.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '431.82'
$ws.Range("E42").Value = '  -3.28%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.287'
$ws.Range("E43").Value = '  +0.51%  '

$ws.Range("E44").Value = '  +2.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0360'
$ws.Range("E45").Value = '  -0.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.827.40'
$ws.Range("E46").Value = '  +0.85%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.23'
$ws.Range("E47").Value = '  -4.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '128.51'
$ws.Range("E48").Value = '  -2.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.97'
$ws.Range("E50").Value = '  -1.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.109'
$ws.Range("E51").Value = '  -1.66%  '
